$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.162.58'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '3.369.89'
$ws.Range('E3').Value = '  +1.79%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'572.03"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('D6').Value = "'137.32"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.83%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.367.61'
$ws.Range('E8').Value = '  +1.81%  '
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('E10').Value = '  +3.73%  '
$ws.Range('E11').Value = '  +4.33%  '
$ws.Range('E12').Value = '  +4.45%  '
$ws.Range('D13').Value = '3.947.49'
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('E14').Value = '  +2.05%  '
$ws.Range('E15').Value = '  +3.54%  '
$ws.Range('D16').Value = '3.368.29'
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('D17').Value = "'25.16"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.77%  '
$ws.Range('D18').Value = '61.263.47'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = "'13.92"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.04%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = "'5.87"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.88%  '
$ws.Range('D21').Value = "'9.35"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.56%  '
$ws.Range('D22').Value = "'377.88"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.89%  '
$ws.Range('D23').Value = "'0.568"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('D24').Value = '3.505.10'
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('D26').Value = "'70.61"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.98%  '
$ws.Range('E27').Value = '  +11.99%  '
$ws.Range('E28').Value = '  +15.96%  '
$ws.Range('D29').Value = "'7.73"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.97%  '
$ws.Range('D30').Value = "'0.994"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('D31').Value = "'8.22"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.68%  '
$ws.Range('E32').Value = '  +4.70%  '
$ws.Range('D33').Value = "'2.12"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.05%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '3.401.01'
$ws.Range('E35').Value = '  +1.97%  '
$ws.Range('D36').Value = "'23.39"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.82%  '
$ws.Range('E37').Value = '  +6.80%  '
$ws.Range('D38').Value = "'7.02"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.40%  '
$ws.Range('E39').Value = '  +5.45%  '
$ws.Range('D40').Value = "'160.90"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.49%  '
$ws.Range('D41').Value = "'0.0789"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.98%  '
$ws.Range('D42').Value = "'1.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('E43').Value = '  +10.63%  '
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('D45').Value = "'41.42"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('E47').Value = '  +7.50%  '
$ws.Range('D48').Value = "'23.00"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.18%  '
$ws.Range('E49').Value = '  +3.96%  '
$ws.Range('D50').Value = "'22.86"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.65%  '
$ws.Range('D51').Value = '2.326.68'
$ws.Range('E51').Value = '  +8.04%  '
